# Update the "Correspond Handoff Datetime" / "Correspond Handback DateTime"
# timestamps for the zh-cn and de-de handback rows, as part of (re)generating
# the handback status report.

$wb = $excel.ActiveWorkbook

# zh-cn sheet (row 2): handoff datetime (E2) and handback datetime (H2)
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E2").Value = "2016-03-22 09:08:27"
$wsZhCn.Range("H2").Value = "2016-03-22 09:10:45"

# de-de sheet (row 2): handoff datetime (E2) and handback datetime (H2)
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E2").Value = "2016-03-22 09:08:31"
$wsDeDe.Range("H2").Value = "2016-03-22 09:10:53"
